# Datorama_Creative_Hierarchcial.xlsx
# Commit: "Added Creative Level Campaign and Advertiser Hierarchy"
#
# Adds two new worksheets ("CreativeDelivery_AdvertiserID" and
# "CreativeDelivery_CampaignID"), each cloned from the first sheet's
# layout/columns, with their own sourceColumn/destinationTableColumnName
# pair written into row 2. Also tweaks the view/selection state on the
# CreativeID sheet and leaves the new CampaignID sheet as the focused tab.

$wb = $excel.ActiveWorkbook

$campaignTargetSheet = $wb.Worksheets.Item(1)
$creativeIdSheet     = $wb.Worksheets.Item(2)

# --- New sheet: CreativeDelivery_AdvertiserID ------------------------------
# Clone the CampaignTarget sheet (it already carries all 5 full-width
# "bestFit" columns) and drop the copy after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$campaignTargetSheet.Copy($null, $lastSheet)
$advertiserIdSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$advertiserIdSheet.Name = "CreativeDelivery_AdvertiserID"
$advertiserIdSheet.Range("A2").Value = "advertiser_id"
$advertiserIdSheet.Range("B2").Value = "Advertiser ID"

# --- New sheet: CreativeDelivery_CampaignID --------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$campaignTargetSheet.Copy($null, $lastSheet)
$campaignIdSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$campaignIdSheet.Name = "CreativeDelivery_CampaignID"
$campaignIdSheet.Range("A2").Value = "campaign_id"
$campaignIdSheet.Range("B2").Value = "Campaign ID"

# --- Selection / view tweaks ------------------------------------------------
# CreativeID sheet: select the whole used range and give it a portrait
# page setup; it is no longer the tab in focus.
[void]$creativeIdSheet.Range("A1:E5").Select()
$creativeIdSheet.PageSetup.Orientation = 1

# AdvertiserID sheet: whole-range selection too.
[void]$advertiserIdSheet.Range("A1:E5").Select()

# CampaignID sheet becomes the active/focused tab, with C2 selected last.
[void]$campaignIdSheet.Range("C2").Select()
